# Fixing name of Sectors to be alligned with Baseline
#
# The workbook contains one sheet per year (2000-2100). Every sheet shares
# the same small header block in row 3, columns D:G, holding the sector
# abbreviations "Nd", "Dy", "Cu", "Si" (each backed by a single shared
# string reused on every sheet). Rename them to their full/aligned names.

$wb = $excel.ActiveWorkbook

$map = @{
    "Nd" = "Neodymium"
    "Dy" = "Dysprosium"
    "Cu" = "Copper ores and concentrates"
    "Si" = "Raw silicon"
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $d3 = $ws.Range("D3").Value()
    if ($map.ContainsKey($d3)) {
        $ws.Range("D3").Value = $map[$d3]
    }

    $e3 = $ws.Range("E3").Value()
    if ($map.ContainsKey($e3)) {
        $ws.Range("E3").Value = $map[$e3]
    }

    $f3 = $ws.Range("F3").Value()
    if ($map.ContainsKey($f3)) {
        $ws.Range("F3").Value = $map[$f3]
    }

    $g3 = $ws.Range("G3").Value()
    if ($map.ContainsKey($g3)) {
        $ws.Range("G3").Value = $map[$g3]
    }
}
